$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data as scraped by the GitHub Actions job.
# D (Price) and E (Volume 1h) columns hold text-formatted numbers/percentages; the
# NumberFormat="@" + ClearFormats() pattern forces Excel to store them as text instead of
# re-parsing them into numeric values, while leaving cell styling untouched.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

# Row 2
Set-TextValue $ws.Cells.Item(2, 4) "63.761.23"
Set-TextValue $ws.Cells.Item(2, 5) "  -1.10%  "
# Row 3
Set-TextValue $ws.Cells.Item(3, 4) "3.077.83"
Set-TextValue $ws.Cells.Item(3, 5) "  -7.67%  "
# Row 4
Set-TextValue $ws.Cells.Item(4, 5) "  -1.33%  "
# Row 5
Set-TextValue $ws.Cells.Item(5, 4) "586.99"
Set-TextValue $ws.Cells.Item(5, 5) "  -1.29%  "
# Row 6
Set-TextValue $ws.Cells.Item(6, 4) "155.65"
Set-TextValue $ws.Cells.Item(6, 5) "  +4.71%  "
# Row 7
Set-TextValue $ws.Cells.Item(7, 4) "0.999"
Set-TextValue $ws.Cells.Item(7, 5) "  -1.36%  "
# Row 8
Set-TextValue $ws.Cells.Item(8, 5) "  +0.22%  "
# Row 9
Set-TextValue $ws.Cells.Item(9, 4) "3.077.34"
Set-TextValue $ws.Cells.Item(9, 5) "  -3.24%  "
# Row 10
Set-TextValue $ws.Cells.Item(10, 5) "  -4.62%  "
# Row 11
Set-TextValue $ws.Cells.Item(11, 5) "  -3.39%  "
# Row 12
Set-TextValue $ws.Cells.Item(12, 4) "0.449"
Set-TextValue $ws.Cells.Item(12, 5) "  -3.59%  "
# Row 13
Set-TextValue $ws.Cells.Item(13, 4) "0.0000237"
Set-TextValue $ws.Cells.Item(13, 5) "  -4.79%  "
# Row 14
Set-TextValue $ws.Cells.Item(14, 4) "36.77"
Set-TextValue $ws.Cells.Item(14, 5) "  -3.43%  "
# Row 15
Set-TextValue $ws.Cells.Item(15, 4) "0.120"
Set-TextValue $ws.Cells.Item(15, 5) "  -1.92%  "
# Row 16
Set-TextValue $ws.Cells.Item(16, 4) "3.584.33"
Set-TextValue $ws.Cells.Item(16, 5) "  -7.63%  "
# Row 17
Set-TextValue $ws.Cells.Item(17, 4) "7.16"
Set-TextValue $ws.Cells.Item(17, 5) "  -2.82%  "
# Row 18
Set-TextValue $ws.Cells.Item(18, 4) "63.668.62"
Set-TextValue $ws.Cells.Item(18, 5) "  -1.15%  "
# Row 19
Set-TextValue $ws.Cells.Item(19, 4) "3.076.49"
Set-TextValue $ws.Cells.Item(19, 5) "  -4.50%  "
# Row 20
Set-TextValue $ws.Cells.Item(20, 4) "470.24"
Set-TextValue $ws.Cells.Item(20, 5) "  -1.29%  "
# Row 21
Set-TextValue $ws.Cells.Item(21, 5) "  -2.96%  "
# Row 22
Set-TextValue $ws.Cells.Item(22, 4) "0.704"
Set-TextValue $ws.Cells.Item(22, 5) "  -6.61%  "
# Row 23
Set-TextValue $ws.Cells.Item(23, 4) "7.51"
Set-TextValue $ws.Cells.Item(23, 5) "  -3.45%  "
# Row 24
Set-TextValue $ws.Cells.Item(24, 4) "2.42"
Set-TextValue $ws.Cells.Item(24, 5) "  -3.17%  "
# Row 25
Set-TextValue $ws.Cells.Item(25, 4) "12.84"
Set-TextValue $ws.Cells.Item(25, 5) "  -6.10%  "
# Row 26
Set-TextValue $ws.Cells.Item(26, 4) "80.39"
Set-TextValue $ws.Cells.Item(26, 5) "  -3.71%  "
# Row 27
Set-TextValue $ws.Cells.Item(27, 4) "10.33"
Set-TextValue $ws.Cells.Item(27, 5) "  +2.70%  "
# Row 28
Set-TextValue $ws.Cells.Item(28, 5) "  -0.30%  "
# Row 29
Set-TextValue $ws.Cells.Item(29, 4) "7.39"
Set-TextValue $ws.Cells.Item(29, 5) "  +1.00%  "
# Row 30
Set-TextValue $ws.Cells.Item(30, 2) "FirstDigitalUSD"
Set-TextValue $ws.Cells.Item(30, 3) "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Cells.Item(30, 4) "1.00"
Set-TextValue $ws.Cells.Item(30, 5) "  -1.10%  "
# Row 31
Set-TextValue $ws.Cells.Item(31, 2) "PancakeSwap"
Set-TextValue $ws.Cells.Item(31, 3) "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Cells.Item(31, 4) "2.66"
Set-TextValue $ws.Cells.Item(31, 5) "  -3.43%  "
# Row 32
Set-TextValue $ws.Cells.Item(32, 4) "2.13"
Set-TextValue $ws.Cells.Item(32, 5) "  -6.60%  "
# Row 33
Set-TextValue $ws.Cells.Item(33, 5) "  -7.81%  "
# Row 34
Set-TextValue $ws.Cells.Item(34, 4) "27.10"
Set-TextValue $ws.Cells.Item(34, 5) "  -4.97%  "
# Row 35
Set-TextValue $ws.Cells.Item(35, 4) "0.0₃0827"
Set-TextValue $ws.Cells.Item(35, 5) "  -3.91%  "
# Row 36
Set-TextValue $ws.Cells.Item(36, 5) "  -3.10%  "
# Row 37
Set-TextValue $ws.Cells.Item(37, 4) "5.97"
Set-TextValue $ws.Cells.Item(37, 5) "  -5.00%  "
# Row 38
Set-TextValue $ws.Cells.Item(38, 4) "3.25"
Set-TextValue $ws.Cells.Item(38, 5) "  -2.53%  "
# Row 39
Set-TextValue $ws.Cells.Item(39, 5) "  -5.68%  "
# Row 40
Set-TextValue $ws.Cells.Item(40, 4) "50.60"
Set-TextValue $ws.Cells.Item(40, 5) "  -2.50%  "
# Row 41
Set-TextValue $ws.Cells.Item(41, 4) "9.13"
Set-TextValue $ws.Cells.Item(41, 5) "  -3.41%  "
# Row 42
Set-TextValue $ws.Cells.Item(42, 4) "431.36"
Set-TextValue $ws.Cells.Item(42, 5) "  -8.53%  "
# Row 43
Set-TextValue $ws.Cells.Item(43, 4) "0.288"
Set-TextValue $ws.Cells.Item(43, 5) "  -4.30%  "
# Row 44
Set-TextValue $ws.Cells.Item(44, 5) "  +0.48%  "
# Row 45
Set-TextValue $ws.Cells.Item(45, 5) "  -4.68%  "
# Row 46
Set-TextValue $ws.Cells.Item(46, 4) "39.84"
Set-TextValue $ws.Cells.Item(46, 5) "  +1.05%  "
# Row 47
Set-TextValue $ws.Cells.Item(47, 4) "2.811.09"
# Row 48
Set-TextValue $ws.Cells.Item(48, 4) "129.90"
Set-TextValue $ws.Cells.Item(48, 5) "  -2.13%  "
# Row 49
Set-TextValue $ws.Cells.Item(49, 4) "0.999"
Set-TextValue $ws.Cells.Item(49, 5) "  +0.03%  "
# Row 50
Set-TextValue $ws.Cells.Item(50, 4) "24.89"
Set-TextValue $ws.Cells.Item(50, 5) "  -0.82%  "
# Row 51
Set-TextValue $ws.Cells.Item(51, 4) "2.21"
Set-TextValue $ws.Cells.Item(51, 5) "  -4.09%  "
